$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 31.90834366666667
$ws.Range("H2").Value = 95.725031
$ws.Range("I2").Value = 0.1125536485145784
$ws.Range("J2").Value = 0.1157863270269485
$ws.Range("M2").Value = 17.96491
$ws.Range("N2").Value = 53.89473
$ws.Range("O2").Value = 0.1605217090392217
$ws.Range("P2").Value = 0.1740895558373204
$ws.Range("Q2").Value = 573.2305222207367
$ws.Range("R2").Value = 5159.074699986631
$ws.Range("S2").Value = 0.01806730401815999
$ws.Range("T2").Value = 0.0201571902441562
$ws.Range("G3").Value = 31.90834366666667
$ws.Range("H3").Value = 95.725031
$ws.Range("I3").Value = 0.1125536485145784
$ws.Range("J3").Value = 0.1157863270269485
$ws.Range("O3").Value = 0.3767843730745145
$ws.Range("P3").Value = 0.4086314838509364
$ws.Range("Q3").Value = 1345.514598834376
$ws.Range("R3").Value = 12109.63138950938
$ws.Range("S3").Value = 0.04240845589281469
$ws.Range("T3").Value = 0.04731393862267178
$ws.Range("G4").Value = 31.90834366666667
$ws.Range("H4").Value = 95.725031
$ws.Range("I4").Value = 0.1125536485145784
$ws.Range("J4").Value = 0.1157863270269485
$ws.Range("M4").Value = 7.700863999999999
$ws.Range("N4").Value = 23.102592
$ws.Range("O4").Value = 0.06880946524967933
$ws.Range("P4").Value = 0.07462547785230264
$ws.Range("Q4").Value = 245.7218150422613
$ws.Range("R4").Value = 2211.496335380352
$ws.Range("S4").Value = 0.007744756366188506
$ws.Range("T4").Value = 0.008640609983149019
$ws.Range("G5").Value = 31.90834366666667
$ws.Range("H5").Value = 95.725031
$ws.Range("I5").Value = 0.1125536485145784
$ws.Range("J5").Value = 0.1157863270269485
$ws.Range("M5").Value = 26.1668075
$ws.Range("N5").Value = 52.33361499999999
$ws.Range("O5").Value = 0.2338080547022124
$ws.Range("P5").Value = 0.1690468769527434
$ws.Range("Q5").Value = 834.9394863695107
$ws.Range("R5").Value = 5009.636918217065
$ws.Range("S5").Value = 0.02631594960883014
$ws.Range("T5").Value = 0.01957331697773467
$ws.Range("G6").Value = 31.90834366666667
$ws.Range("H6").Value = 95.725031
$ws.Range("I6").Value = 0.1125536485145784
$ws.Range("J6").Value = 0.1157863270269485
$ws.Range("M6").Value = 17.91507266666667
$ws.Range("N6").Value = 53.74521800000001
$ws.Range("O6").Value = 0.1600763979343721
$ws.Range("P6").Value = 0.1736066055066972
$ws.Range("Q6").Value = 571.6402954613066
$ws.Range("R6").Value = 5144.762659151759
$ws.Range("S6").Value = 0.01801718262858511
$ws.Range("T6").Value = 0.02010127119923689
$ws.Range("I7").Value = 0.2312918537506949
$ws.Range("J7").Value = 0.2379348388122522
$ws.Range("M7").Value = 17.96491
$ws.Range("N7").Value = 53.89473
$ws.Range("O7").Value = 0.1605217090392217
$ws.Range("P7").Value = 0.1740895558373204
$ws.Range("Q7").Value = 1177.958705565553
$ws.Range("R7").Value = 10601.62835008998
$ws.Range("S7").Value = 0.03712736365091126
$ws.Range("T7").Value = 0.04142197040704941
$ws.Range("I8").Value = 0.2312918537506949
$ws.Range("J8").Value = 0.2379348388122522
$ws.Range("O8").Value = 0.3767843730745145
$ws.Range("P8").Value = 0.4086314838509364
$ws.Range("S8").Value = 0.08714715611269785
$ws.Range("T8").Value = 0.09722766624368402
$ws.Range("I9").Value = 0.2312918537506949
$ws.Range("J9").Value = 0.2379348388122522
$ws.Range("M9").Value = 7.700863999999999
$ws.Range("N9").Value = 23.102592
$ws.Range("O9").Value = 0.06880946524967933
$ws.Range("P9").Value = 0.07462547785230264
$ws.Range("Q9").Value = 504.9454625253546
$ws.Range("R9").Value = 4544.509162728192
$ws.Range("S9").Value = 0.01591506877319235
$ws.Range("T9").Value = 0.01775600104407493
$ws.Range("I10").Value = 0.2312918537506949
$ws.Range("J10").Value = 0.2379348388122522
$ws.Range("M10").Value = 26.1668075
$ws.Range("N10").Value = 52.33361499999999
$ws.Range("O10").Value = 0.2338080547022124
$ws.Range("P10").Value = 0.1690468769527434
$ws.Range("Q10").Value = 1715.756922327081
$ws.Range("R10").Value = 10294.54153396249
$ws.Range("S10").Value = 0.05407789839391858
$ws.Range("T10").Value = 0.04022214141946563
$ws.Range("I11").Value = 0.2312918537506949
$ws.Range("J11").Value = 0.2379348388122522
$ws.Range("M11").Value = 17.91507266666667
$ws.Range("N11").Value = 53.74521800000001
$ws.Range("O11").Value = 0.1600763979343721
$ws.Range("P11").Value = 0.1736066055066972
$ws.Range("Q11").Value = 1174.690872848208
$ws.Range("R11").Value = 10572.21785563387
$ws.Range("S11").Value = 0.03702436681997483
$ws.Range("T11").Value = 0.04130705969797826
$ws.Range("G12").Value = 85.57939900000001
$ws.Range("H12").Value = 256.738197
$ws.Range("I12").Value = 0.3018731932863474
$ws.Range("J12").Value = 0.3105433607867011
$ws.Range("M12").Value = 17.96491
$ws.Range("N12").Value = 53.89473
$ws.Range("O12").Value = 0.1605217090392217
$ws.Range("P12").Value = 0.1740895558373204
$ws.Range("Q12").Value = 1537.42620088909
$ws.Range("R12").Value = 13836.83580800181
$ws.Range("S12").Value = 0.04845720089945178
$ws.Range("T12").Value = 0.05406235574758553
$ws.Range("G13").Value = 85.57939900000001
$ws.Range("H13").Value = 256.738197
$ws.Range("I13").Value = 0.3018731932863474
$ws.Range("J13").Value = 0.3105433607867011
$ws.Range("O13").Value = 0.3767843730745145
$ws.Range("P13").Value = 0.4086314838509364
$ws.Range("Q13").Value = 3608.721653398222
$ws.Range("R13").Value = 32478.494880584
$ws.Range("S13").Value = 0.1137411018803981
$ws.Range("T13").Value = 0.1268977943183264
$ws.Range("G14").Value = 85.57939900000001
$ws.Range("H14").Value = 256.738197
$ws.Range("I14").Value = 0.3018731932863474
$ws.Range("J14").Value = 0.3105433607867011
$ws.Range("M14").Value = 7.700863999999999
$ws.Range("N14").Value = 23.102592
$ws.Range("O14").Value = 0.06880946524967933
$ws.Range("P14").Value = 0.07462547785230264
$ws.Range("Q14").Value = 659.035312900736
$ws.Range("R14").Value = 5931.317816106623
$ws.Range("S14").Value = 0.02077173300324665
$ws.Range("T14").Value = 0.02317444669256759
$ws.Range("G15").Value = 85.57939900000001
$ws.Range("H15").Value = 256.738197
$ws.Range("I15").Value = 0.3018731932863474
$ws.Range("J15").Value = 0.3105433607867011
$ws.Range("M15").Value = 26.1668075
$ws.Range("N15").Value = 52.33361499999999
$ws.Range("O15").Value = 0.2338080547022124
$ws.Range("P15").Value = 0.1690468769527434
$ws.Range("Q15").Value = 2239.339659598692
$ws.Range("R15").Value = 13436.03795759215
$ws.Range("S15").Value = 0.07058038408902585
$ws.Range("T15").Value = 0.05249638529940085
$ws.Range("G16").Value = 85.57939900000001
$ws.Range("H16").Value = 256.738197
$ws.Range("I16").Value = 0.3018731932863474
$ws.Range("J16").Value = 0.3105433607867011
$ws.Range("M16").Value = 17.91507266666667
$ws.Range("N16").Value = 53.74521800000001
$ws.Range("O16").Value = 0.1600763979343721
$ws.Range("P16").Value = 0.1736066055066972
$ws.Range("Q16").Value = 1533.161151854661
$ws.Range("R16").Value = 13798.45036669195
$ws.Range("S16").Value = 0.04832277341422497
$ws.Range("T16").Value = 0.05391237872882075
$ws.Range("G17").Value = 23.7449455
$ws.Range("H17").Value = 47.489891
$ws.Range("I17").Value = 0.08375803763818537
$ws.Range("J17").Value = 0.05744244731349463
$ws.Range("M17").Value = 17.96491
$ws.Range("N17").Value = 53.89473
$ws.Range("O17").Value = 0.1605217090392217
$ws.Range("P17").Value = 0.1740895558373204
$ws.Range("Q17").Value = 426.575808862405
$ws.Range("R17").Value = 2559.45485317443
$ws.Range("S17").Value = 0.01344498334745297
$ws.Range("T17").Value = 0.01000013013901496
$ws.Range("G18").Value = 23.7449455
$ws.Range("H18").Value = 47.489891
$ws.Range("I18").Value = 0.08375803763818537
$ws.Range("J18").Value = 0.05744244731349463
$ws.Range("O18").Value = 0.3767843730745145
$ws.Range("P18").Value = 0.4086314838509364
$ws.Range("Q18").Value = 1001.279513362914
$ws.Range("R18").Value = 6007.677080177485
$ws.Range("S18").Value = 0.03155871970145527
$ws.Range("T18").Value = 0.02347279248174255
$ws.Range("G19").Value = 23.7449455
$ws.Range("H19").Value = 47.489891
$ws.Range("I19").Value = 0.08375803763818537
$ws.Range("J19").Value = 0.05744244731349463
$ws.Range("M19").Value = 7.700863999999999
$ws.Range("N19").Value = 23.102592
$ws.Range("O19").Value = 0.06880946524967933
$ws.Range("P19").Value = 0.07462547785230264
$ws.Range("Q19").Value = 182.856595982912
$ws.Range("R19").Value = 1097.139575897472
$ws.Range("S19").Value = 0.00576334578024605
$ws.Range("T19").Value = 0.004286670079775255
$ws.Range("G20").Value = 23.7449455
$ws.Range("H20").Value = 47.489891
$ws.Range("I20").Value = 0.08375803763818537
$ws.Range("J20").Value = 0.05744244731349463
$ws.Range("M20").Value = 26.1668075
$ws.Range("N20").Value = 52.33361499999999
$ws.Range("O20").Value = 0.2338080547022124
$ws.Range("P20").Value = 0.1690468769527434
$ws.Range("Q20").Value = 621.3294179964912
$ws.Range("R20").Value = 2485.317671985965
$ws.Range("S20").Value = 0.01958330384585881
$ws.Range("T20").Value = 0.009710466322868771
$ws.Range("G21").Value = 23.7449455
$ws.Range("H21").Value = 47.489891
$ws.Range("I21").Value = 0.08375803763818537
$ws.Range("J21").Value = 0.05744244731349463
$ws.Range("M21").Value = 17.91507266666667
$ws.Range("N21").Value = 53.74521800000001
$ws.Range("O21").Value = 0.1600763979343721
$ws.Range("P21").Value = 0.1736066055066972
$ws.Range("Q21").Value = 425.3924240985398
$ws.Range("R21").Value = 2552.354544591238
$ws.Range("S21").Value = 0.01340768496317228
$ws.Range("T21").Value = 0.009972388290093101
$ws.Range("G22").Value = 76.69186633333334
$ws.Range("H22").Value = 230.075599
$ws.Range("I22").Value = 0.270523266810194
$ws.Range("J22").Value = 0.2782930260606035
$ws.Range("M22").Value = 17.96491
$ws.Range("N22").Value = 53.89473
$ws.Range("O22").Value = 0.1605217090392217
$ws.Range("P22").Value = 0.1740895558373204
$ws.Range("Q22").Value = 1377.762476410363
$ws.Range("R22").Value = 12399.86228769327
$ws.Range("S22").Value = 0.0434248571232457
$ws.Range("T22").Value = 0.04844790929951429
$ws.Range("G23").Value = 76.69186633333334
$ws.Range("H23").Value = 230.075599
$ws.Range("I23").Value = 0.270523266810194
$ws.Range("J23").Value = 0.2782930260606035
$ws.Range("O23").Value = 0.3767843730745145
$ws.Range("P23").Value = 0.4086314838509364
$ws.Range("Q23").Value = 3233.951183469074
$ws.Range("R23").Value = 29105.56065122167
$ws.Range("S23").Value = 0.1019289394871486
$ws.Range("T23").Value = 0.1137192921845117
$ws.Range("G24").Value = 76.69186633333334
$ws.Range("H24").Value = 230.075599
$ws.Range("I24").Value = 0.270523266810194
$ws.Range("J24").Value = 0.2782930260606035
$ws.Range("M24").Value = 7.700863999999999
$ws.Range("N24").Value = 23.102592
$ws.Range("O24").Value = 0.06880946524967933
$ws.Range("P24").Value = 0.07462547785230264
$ws.Range("Q24").Value = 590.5936325391787
$ws.Range("R24").Value = 5315.342692852608
$ws.Range("S24").Value = 0.01861456132680578
$ws.Range("T24").Value = 0.02076775005273585
$ws.Range("G25").Value = 76.69186633333334
$ws.Range("H25").Value = 230.075599
$ws.Range("I25").Value = 0.270523266810194
$ws.Range("J25").Value = 0.2782930260606035
$ws.Range("M25").Value = 26.1668075
$ws.Range("N25").Value = 52.33361499999999
$ws.Range("O25").Value = 0.2338080547022124
$ws.Range("P25").Value = 0.1690468769527434
$ws.Range("Q25").Value = 2006.781303160064
$ws.Range("R25").Value = 12040.68781896038
$ws.Range("S25").Value = 0.06325051876457906
$ws.Range("T25").Value = 0.04704456693327345
$ws.Range("G26").Value = 76.69186633333334
$ws.Range("H26").Value = 230.075599
$ws.Range("I26").Value = 0.270523266810194
$ws.Range("J26").Value = 0.2782930260606035
$ws.Range("M26").Value = 17.91507266666667
$ws.Range("N26").Value = 53.74521800000001
$ws.Range("O26").Value = 0.1600763979343721
$ws.Range("P26").Value = 0.1736066055066972
$ws.Range("Q26").Value = 1373.940358303954
$ws.Range("R26").Value = 12365.46322473558
$ws.Range("S26").Value = 0.04330439010841493
$ws.Range("T26").Value = 0.04831350759056819
